$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.655228775655152
$ws.Range("D2").Value = 0.5191137396196956

# Row 3
$ws.Range("C3").Value = 1.833536124438842
$ws.Range("D3").Value = 0.08028864824661763
$ws.Range("G3").Value = "No"

# Row 4
$ws.Range("C4").Value = 0.8650511885208452
$ws.Range("D4").Value = 0.3963457664931593

# Row 5
$ws.Range("C5").Value = 0.5942906078244089
$ws.Range("D5").Value = 0.5583803462765995

# Row 6
$ws.Range("C6").Value = 1.649164664024464
$ws.Range("D6").Value = 0.1133214844882795

# Row 7
$ws.Range("C7").Value = 0.1176955093445544
$ws.Range("D7").Value = 0.9073768684073185

# Row 8
$ws.Range("C8").Value = 0.06159985544699942
$ws.Range("D8").Value = 0.9514377015822293

# Row 9
$ws.Range("C9").Value = -1.147397209234486
$ws.Range("D9").Value = 0.2635426811050379

# Row 10
$ws.Range("C10").Value = -0.780905884131006
$ws.Range("D10").Value = 0.4431802023019524

# Row 11
$ws.Range("C11").Value = -0.0217580378119206
$ws.Range("D11").Value = 0.9828371069407114
